$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Percent Complete (E) and Has Description (G) values for rooms that
# gained descriptions/content.
$ws.Range("E4").Value = 10
$ws.Range("G4").Value = 1

$ws.Range("E5").Value = 10
$ws.Range("G5").Value = 1

$ws.Range("E6").Value = 20
$ws.Range("G6").Value = 1

$ws.Range("E10").Value = 20

$ws.Range("E12").Value = 10
$ws.Range("G12").Value = 1

# Move the active selection to G6 to reflect where editing last occurred.
$ws.Range("G6").Select()
